$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B31").Value = 16
$ws.Range("D31").Value = 88
$ws.Range("G31").Value = 6.571428571428571

$ws.Range("B33").Value = 2
